$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "G5"
$ws.Range("B6").Value = "Investment Plan"
$ws.Range("C6").Value = "Daily"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 45907
$ws.Range("E6").NumberFormat = $ws.Range("E5").NumberFormat
$ws.Range("F6").Value = 36
